$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("General")

# Copy the formatting of the previous data row (title + wrapped description)
# onto the new row 11 so the new cells inherit the same cell styles
# (bold/underline title style + bordered, wrapped description style).
$ws.Range("B10:C10").Copy()
$ws.Range("B11:C11").PasteSpecial(-4122)

$ws.Range("B11").Value = "_C3D-TEMPLATE_2025_FRA (Architecture v0001g)"
$ws.Range("C11").Value = "Renommage des Définitions des jeux de propriétés ACA : 
- ACA-JPPA-objet (pour les jeux personnalisés de propriétés applicables - JPPA) ;
- ACA-Style-objet (pour les jeux de propriétés automatiques, hérités des styles d'objets ACA, notion qui n'existait pas jusqu'alors pour la modélisation traditionnelle avec Civil 3D)"

$ws.Rows.Item(11).RowHeight = 57

$ws.Application.CutCopyMode = $false

$ws.Range("G10").Select()
